$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 596585.0600000001
$ws.Range("I15").Value = 596585.0600000001
$ws.Range("K15").Value = 1789755.18
$ws.Range("M15").Value = -1789586.18
# Row 103
$ws.Range("H103").Value = 1299.6666
$ws.Range("I103").Value = 1299.6666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 3898.9998
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -3312.9998
$ws.Range("N103").ClearContents()
# Row 107
$ws.Range("H107").Value = 915.82355
$ws.Range("I107").Value = 985.7143
$ws.Range("J107").Value = 589.6667
$ws.Range("K107").Value = 985.7143
$ws.Range("L107").Value = 589.6667
$ws.Range("M107").Value = 934.2857
$ws.Range("N107").Value = -4429.6667
# Row 129
$ws.Range("H129").Value = 5862.6665
$ws.Range("I129").Value = 4337.5
$ws.Range("J129").Value = 10743.2
$ws.Range("K129").Value = 13012.5
$ws.Range("L129").Value = 32229.6
$ws.Range("M129").Value = -8012.5
$ws.Range("N129").Value = -42229.60000000001
# Row 138
$ws.Range("H138").Value = 10055.659
$ws.Range("I138").Value = 4999.5
$ws.Range("J138").Value = 10280.378
$ws.Range("K138").Value = 14998.5
$ws.Range("L138").Value = 30841.134
$ws.Range("M138").Value = -9858.5
$ws.Range("N138").Value = -41121.13400000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5648.373
$ws.Range("I32").Value = 4335.295
$ws.Range("J32").Value = 18998
$ws.Range("K32").Value = 4335.295
$ws.Range("L32").Value = 18998
$ws.Range("M32").Value = -4048.295
$ws.Range("N32").Value = -19572
# Row 45
$ws.Range("H45").Value = 2083.55
$ws.Range("I45").Value = 1408.7858
$ws.Range("J45").Value = 3658
$ws.Range("K45").Value = 1408.7858
$ws.Range("L45").Value = 3658
$ws.Range("M45").Value = -1031.7858
$ws.Range("N45").Value = -4412
# Row 61
$ws.Range("H61").Value = 4540.5
$ws.Range("I61").Value = 4540.5
$ws.Range("K61").Value = 4540.5
$ws.Range("M61").Value = -4328.5
# Row 74
$ws.Range("H74").Value = 8445.85
$ws.Range("I74").Value = 3946.7693
$ws.Range("J74").Value = 16801.285
$ws.Range("K74").Value = 3946.7693
$ws.Range("L74").Value = 16801.285
$ws.Range("M74").Value = -3072.7693
$ws.Range("N74").Value = -18549.285
# Row 77
$ws.Range("H77").Value = 8445.85
$ws.Range("I77").Value = 3946.7693
$ws.Range("J77").Value = 16801.285
$ws.Range("K77").Value = 19733.8465
$ws.Range("L77").Value = 84006.425
$ws.Range("M77").Value = -15365.8465
$ws.Range("N77").Value = -92742.425
# Row 102
$ws.Range("H102").Value = 1508.0714
$ws.Range("I102").Value = 854.9231
$ws.Range("J102").Value = 9999
$ws.Range("K102").Value = 854.9231
$ws.Range("L102").Value = 9999
$ws.Range("M102").Value = 767.0769
$ws.Range("N102").Value = -13243
# Row 131
$ws.Range("H131").Value = 41238.332
$ws.Range("I131").Value = 50000
$ws.Range("J131").Value = 36857.5
$ws.Range("K131").Value = 50000
$ws.Range("L131").Value = 36857.5
$ws.Range("M131").Value = -44960
$ws.Range("N131").Value = -46937.5
# Row 132
$ws.Range("H132").Value = 9807.454
$ws.Range("I132").Value = 6563.8335
$ws.Range("J132").Value = 13699.8
$ws.Range("K132").Value = 19691.5005
$ws.Range("L132").Value = 41099.39999999999
$ws.Range("M132").Value = -17161.5005
$ws.Range("N132").Value = -46159.39999999999
# Row 136
$ws.Range("H136").Value = 4540.5
$ws.Range("I136").Value = 4540.5
$ws.Range("K136").Value = 13621.5
$ws.Range("M136").Value = -11071.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3923.762
$ws.Range("I86").Value = 1682.091
$ws.Range("J86").Value = 6389.6
$ws.Range("K86").Value = 1682.091
$ws.Range("L86").Value = 6389.6
$ws.Range("M86").Value = -559.0909999999999
$ws.Range("N86").Value = -8635.6
# Row 89
$ws.Range("H89").Value = 3923.762
$ws.Range("I89").Value = 1682.091
$ws.Range("J89").Value = 6389.6
$ws.Range("K89").Value = 8410.455
$ws.Range("L89").Value = 31948
$ws.Range("M89").Value = -2794.455
$ws.Range("N89").Value = -43180
# Row 94
$ws.Range("H94").Value = 882.4
$ws.Range("I94").Value = 561.6842
$ws.Range("J94").Value = 1898
$ws.Range("K94").Value = 561.6842
$ws.Range("L94").Value = 1898
$ws.Range("M94").Value = -110.6842
$ws.Range("N94").Value = -2800
# Row 105
$ws.Range("H105").Value = 2703.037
$ws.Range("I105").Value = 2412.8333
$ws.Range("J105").Value = 3283.4443
$ws.Range("K105").Value = 2412.8333
$ws.Range("L105").Value = 3283.4443
$ws.Range("M105").Value = -665.8332999999998
$ws.Range("N105").Value = -6777.4443
# Row 134
$ws.Range("H134").Value = 3015.75
$ws.Range("I134").Value = 1348.1765
$ws.Range("J134").Value = 12465.333
$ws.Range("K134").Value = 4044.5295
$ws.Range("L134").Value = 37395.999
$ws.Range("M134").Value = -1509.5295
$ws.Range("N134").Value = -42465.999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1955.2222
$ws.Range("I16").Value = 949.5
$ws.Range("J16").Value = 3966.6667
$ws.Range("K16").Value = 949.5
$ws.Range("L16").Value = 3966.6667
$ws.Range("M16").Value = -662.5
$ws.Range("N16").Value = -4540.6667
# Row 31
$ws.Range("H31").Value = 5038.956
$ws.Range("I31").Value = 1827.5625
$ws.Range("J31").Value = 7893.528
$ws.Range("K31").Value = 1827.5625
$ws.Range("L31").Value = 7893.528
$ws.Range("M31").Value = -1532.5625
$ws.Range("N31").Value = -8483.528
# Row 34
$ws.Range("H34").Value = 5038.956
$ws.Range("I34").Value = 1827.5625
$ws.Range("J34").Value = 7893.528
$ws.Range("K34").Value = 1827.5625
$ws.Range("L34").Value = 7893.528
$ws.Range("M34").Value = -1625.5625
$ws.Range("N34").Value = -8297.528
# Row 105
$ws.Range("H105").Value = 1282.625
$ws.Range("I105").Value = 1251.1
$ws.Range("K105").Value = 1251.1
$ws.Range("M105").Value = 495.9000000000001
# Row 107
$ws.Range("H107").Value = 530.2258
$ws.Range("I107").Value = 337.68182
$ws.Range("J107").Value = 1000.8889
$ws.Range("K107").Value = 337.68182
$ws.Range("L107").Value = 1000.8889
$ws.Range("M107").Value = 1582.31818
$ws.Range("N107").Value = -4840.8889
# Row 113
$ws.Range("H113").Value = 1955.2222
$ws.Range("I113").Value = 949.5
$ws.Range("J113").Value = 3966.6667
$ws.Range("K113").Value = 949.5
$ws.Range("L113").Value = 3966.6667
$ws.Range("M113").Value = 1220.5
$ws.Range("N113").Value = -8306.6667
# Row 122
$ws.Range("H122").Value = 1743.2667
$ws.Range("I122").Value = 1479.0834
$ws.Range("K122").Value = 4437.2502
$ws.Range("M122").Value = -1987.2502
# Row 132
$ws.Range("H132").Value = 5149.9
$ws.Range("I132").Value = 3374.75
$ws.Range("K132").Value = 10124.25
$ws.Range("M132").Value = -7594.25
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 4748.174
$ws.Range("J68").Value = 7201.2
$ws.Range("L68").Value = 21603.6
$ws.Range("N68").Value = -23225.6
# Row 71
$ws.Range("H71").Value = 4748.174
$ws.Range("J71").Value = 7201.2
$ws.Range("L71").Value = 64810.8
$ws.Range("N71").Value = -72922.79999999999
# Row 75
$ws.Range("H75").Value = 13238.875
$ws.Range("I75").Value = 955.5
$ws.Range("J75").Value = 17333.334
$ws.Range("K75").Value = 2866.5
$ws.Range("L75").Value = 52000.00199999999
$ws.Range("M75").Value = -1868.5
$ws.Range("N75").Value = -53996.00199999999
# Row 78
$ws.Range("H78").Value = 13238.875
$ws.Range("I78").Value = 955.5
$ws.Range("J78").Value = 17333.334
$ws.Range("K78").Value = 8599.5
$ws.Range("L78").Value = 156000.006
$ws.Range("M78").Value = -3607.5
$ws.Range("N78").Value = -165984.006
# Row 107
$ws.Range("H107").Value = 2000
$ws.Range("J107").Value = 2000
$ws.Range("L107").Value = 6000
$ws.Range("N107").Value = -9840
# Row 119
$ws.Range("H119").Value = 920
$ws.Range("I119").Value = 920
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 2760
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = 2078
$ws.Range("N119").ClearContents()
# Row 131
$ws.Range("H131").Value = 700473.5600000001
$ws.Range("I131").Value = 863.3158
$ws.Range("J131").Value = 2599415.5
$ws.Range("K131").Value = 2589.9474
$ws.Range("L131").Value = 7798246.5
$ws.Range("M131").Value = 2450.0526
$ws.Range("N131").Value = -7808326.5
# Row 132
$ws.Range("H132").Value = 2725.3157
$ws.Range("I132").Value = 2522.6924
$ws.Range("J132").Value = 3164.3333
$ws.Range("K132").Value = 22704.2316
$ws.Range("L132").Value = 28478.9997
$ws.Range("M132").Value = -20174.2316
$ws.Range("N132").Value = -33538.9997
# Row 138
$ws.Range("H138").Value = 1853.375
$ws.Range("I138").Value = 1853.375
$ws.Range("K138").Value = 5560.125
$ws.Range("M138").Value = -420.125
# Row 140
$ws.Range("H140").Value = 1630.65
$ws.Range("I140").Value = 1514.0526
$ws.Range("J140").Value = 1736.1428
$ws.Range("K140").Value = 4542.1578
$ws.Range("L140").Value = 5208.428400000001
$ws.Range("M140").Value = 637.8422
$ws.Range("N140").Value = -15568.4284

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 11782.182
$ws.Range("J80").Value = 12900.833
$ws.Range("L80").Value = 12900.833
$ws.Range("N80").Value = -14896.833
# Row 83
$ws.Range("H83").Value = 11782.182
$ws.Range("J83").Value = 12900.833
$ws.Range("L83").Value = 64504.165
$ws.Range("N83").Value = -74488.16500000001
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 107
$ws.Range("H107").Value = 904
$ws.Range("I107").Value = 509.8
$ws.Range("J107").Value = 1396.75
$ws.Range("K107").Value = 509.8
$ws.Range("L107").Value = 1396.75
$ws.Range("M107").Value = 1410.2
$ws.Range("N107").Value = -5236.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 53
$ws.Range("H53").Value = 17957
$ws.Range("I53").Value = 3950
$ws.Range("J53").Value = 23559.8
$ws.Range("K53").Value = 3950
$ws.Range("L53").Value = 23559.8
$ws.Range("M53").Value = -3432
$ws.Range("N53").Value = -24595.8
# Row 122
$ws.Range("H122").Value = 8798.647000000001
$ws.Range("I122").Value = 6848
$ws.Range("K122").Value = 20544
$ws.Range("M122").Value = -18094
# Row 132
$ws.Range("H132").Value = 9329.368
$ws.Range("I132").Value = 5061.8
$ws.Range("J132").Value = 14071.111
$ws.Range("K132").Value = 15185.4
$ws.Range("L132").Value = 42213.333
$ws.Range("M132").Value = -12655.4
$ws.Range("N132").Value = -47273.333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3912.5454
$ws.Range("I122").Value = 3152.8462
$ws.Range("K122").Value = 9458.5386
$ws.Range("M122").Value = -7008.5386
# Row 141
$ws.Range("H141").Value = 67999
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 67999
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 67999
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -78359
